# Fix manual mode for money leaderboard
# - Bump the "last updated" timestamp on every leaderboard sheet.
# - Correct the manually-entered tie-break ordering of player names that
#   share the same score (names had been swapped/misplaced between rows).

$wb = $excel.ActiveWorkbook

$newTimestamp = "Dernière update le 05.03.25 à 01:59"

# --- Sheet 1: leaderboard2 ("Qui a attrapé le plus de Cobblemons ?") ---
$ws1 = $wb.Worksheets.Item("leaderboard2")
$ws1.Range("B14").Value = $newTimestamp

# --- Sheet 2: leaderboard3 ("Qui a attrapé le plus de Shiny Cobblemons ?") ---
$ws2 = $wb.Worksheets.Item("leaderboard3")
$ws2.Range("I8").Value  = "LuttiLutti"
$ws2.Range("L8").Value  = "nemenems"
$ws2.Range("I9").Value  = "KyriaaTV"
$ws2.Range("L9").Value  = "Maxouzboub"
$ws2.Range("I10").Value = "Terraciid"
$ws2.Range("B14").Value = $newTimestamp

# --- Sheet 3: leaderboard4 ("Qui a attrapé le plus de Cobblemons légendaires ?") ---
$ws3 = $wb.Worksheets.Item("leaderboard4")
$ws3.Range("I3").Value  = "Kaatsup"
$ws3.Range("O3").Value  = "Elspawn"
$ws3.Range("I4").Value  = "Brybry_"
$ws3.Range("I5").Value  = "Terraciid"
$ws3.Range("O5").Value  = "HarryLafranc"
$ws3.Range("L6").Value  = "ARELIANN"
$ws3.Range("O6").Value  = "Pepito_kawazakii"
$ws3.Range("L7").Value  = "Mynth0s"
$ws3.Range("O7").Value  = "Horty_"
$ws3.Range("L9").Value  = "Maxouzboub"
$ws3.Range("L10").Value = "Angle_Droit"
$ws3.Range("L12").Value = "BagheraJones"
$ws3.Range("F13").Value = "JLTootmy"
$ws3.Range("B14").Value = $newTimestamp

# --- Sheet 4: leaderboard5 ("Qui a le plus de golds ?") ---
$ws4 = $wb.Worksheets.Item("leaderboard5")
$ws4.Range("F3").Value  = "AntoineDaniel_"
$ws4.Range("L3").Value  = "Terraciid"
$ws4.Range("O3").Value  = "Horty_"
$ws4.Range("I4").Value  = "LittleBigWhale"
$ws4.Range("L4").Value  = "Gom4rt_"
$ws4.Range("O4").Value  = "HarryLafranc"
$ws4.Range("I5").Value  = "NakaStream"
$ws4.Range("L5").Value  = "Wingobear"
$ws4.Range("O5").Value  = "Hiro_Ammar"
$ws4.Range("L6").Value  = "ZeratoR"
$ws4.Range("O6").Value  = "CrocodyleTV"
$ws4.Range("I7").Value  = "Bytell2"
$ws4.Range("O7").Value  = "BagheraJones"
$ws4.Range("F8").Value  = "_Linca"
$ws4.Range("I8").Value  = "Angle_Droit"
$ws4.Range("L8").Value  = "JLKada"
$ws4.Range("O8").Value  = "nisqylegoat"
$ws4.Range("F9").Value  = "Etoiles"
$ws4.Range("L9").Value  = "Pepito_kawazakii"
$ws4.Range("F10").Value = "Onutrem"
$ws4.Range("I10").Value = "KennyStream"
$ws4.Range("L10").Value = "XoTrixy"
$ws4.Range("I11").Value = "Kaatsup"
$ws4.Range("L11").Value = "JLTootmy"
$ws4.Range("L12").Value = "Maxouzboub"
$ws4.Range("C13").Value = "ChloeRamdani"
$ws4.Range("I13").Value = "ARELIANN"
$ws4.Range("L13").Value = "Grimkujow"
$ws4.Range("B14").Value = $newTimestamp
